$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") entirely - shifts D->C (prediction) and E->D (rejection-f) left
$ws.Range("C1").EntireColumn.Delete()

# Update the numeric value in B2 (was 1, now 1684.22140908781)
$ws.Range("B2").Value = 1684.22140908781
